$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 26: fill in the previously-empty Reflection cell (F26) ---
$ws.Range("F11").Copy()
$ws.Range("F26").PasteSpecial(-4122)
$ws.Range("F26").Value = "Learnt more about what an expert would do in reverse engineering. Also after learning this course for 5 weeks, we are gonna have a midterm next week. Will review the knowledges and practice tools that we learnt before in the next few days and prepare for the midterm."

# --- Row 27: fill in the previously-empty Reflection cell (F27) ---
$ws.Range("F11").Copy()
$ws.Range("F27").PasteSpecial(-4122)
$ws.Range("F27").Value = "Midterm is kinda hard to describe, since it has many subjective questions. But it’s also the meaning of this course. There is no certain ways to reverse engineer, what we can do is to be subjective and try our best based on some useful concepts. I hope i did it well. Also about the lecture, it is kinda interesting to think about the big picture of one program. Thinking of it stakeholder, functionality and key developers could let us know more than the program itself, like, the community, and the future."

# --- Row 28: new diary entry (was a blank template row) ---
$ws.Range("A27").Copy()
$ws.Range("A28").PasteSpecial(-4122)

$ws.Range("C27").Copy()
$ws.Range("B28:F28").PasteSpecial(-4122)

$ws.Range("G27").Copy()
$ws.Range("G28").PasteSpecial(-4122)

$ws.Range("A28").Value = "2/19/2020"
$ws.Range("B28").Value = "13:00-17:00"
$ws.Range("C28").Value = "Soobin, Marc"
$ws.Range("D28").Value = "Finish homework,3 resubmit homework2"
$ws.Range("E28").Value = "finished homework2, and resubmit our homework3"
$ws.Range("F28").Value = "For our previous homework report, we didn’t realize how the report structure will influence reader’s readability. We change the structure of our homework report and make the content more logical this time."
$ws.Range("G28").Value = "Average"
